$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$wsForecast.Cells.Item(2, 1).Value = 45263.99999999999
$wsForecast.Cells.Item(2, 2).Value = 140
$wsForecast.Cells.Item(2, 3).Value = -91.09666806962306
$wsForecast.Cells.Item(2, 4).Value = 362.8096692810992
$wsForecast.Cells.Item(3, 1).Value = 45270.99999999999
$wsForecast.Cells.Item(3, 2).Value = 142
$wsForecast.Cells.Item(3, 3).Value = -96.30786673262037
$wsForecast.Cells.Item(3, 4).Value = 389.5449260473033
$wsForecast.Cells.Item(4, 1).Value = 45277.99999999999
$wsForecast.Cells.Item(4, 2).Value = 143
$wsForecast.Cells.Item(4, 3).Value = -93.47453745893793
$wsForecast.Cells.Item(4, 4).Value = 373.5474113090217
$wsForecast.Cells.Item(5, 1).Value = 45298.99999999999
$wsForecast.Cells.Item(5, 2).Value = 149
$wsForecast.Cells.Item(5, 3).Value = -77.39052647108484
$wsForecast.Cells.Item(5, 4).Value = 399.3189921207153
$wsForecast.Cells.Item(6, 1).Value = 45305.99999999999
$wsForecast.Cells.Item(6, 2).Value = 151
$wsForecast.Cells.Item(6, 3).Value = -86.61796734404506
$wsForecast.Cells.Item(6, 4).Value = 389.8146292686172
$wsForecast.Cells.Item(7, 1).Value = 45312.99999999999
$wsForecast.Cells.Item(7, 2).Value = 153
$wsForecast.Cells.Item(7, 3).Value = -93.8507210988099
$wsForecast.Cells.Item(7, 4).Value = 385.3369960853288
$wsForecast.Cells.Item(8, 1).Value = 45319.99999999999
$wsForecast.Cells.Item(8, 2).Value = 154
$wsForecast.Cells.Item(8, 3).Value = -84.89189460684409
$wsForecast.Cells.Item(8, 4).Value = 379.1131742034493
$wsForecast.Cells.Item(9, 1).Value = 45326.99999999999
$wsForecast.Cells.Item(9, 2).Value = 156
$wsForecast.Cells.Item(9, 3).Value = -72.07170567281848
$wsForecast.Cells.Item(9, 4).Value = 407.2595616043554
$wsForecast.Cells.Item(10, 1).Value = 45333.99999999999
$wsForecast.Cells.Item(10, 2).Value = 158
$wsForecast.Cells.Item(10, 3).Value = -74.98913819128953
$wsForecast.Cells.Item(10, 4).Value = 391.9818208711765
$wsForecast.Cells.Item(11, 1).Value = 45340.99999999999
$wsForecast.Cells.Item(11, 2).Value = 160
$wsForecast.Cells.Item(11, 3).Value = -76.81006925636945
$wsForecast.Cells.Item(11, 4).Value = 389.3381108529385
$wsForecast.Cells.Item(12, 1).Value = 45347.99999999999
$wsForecast.Cells.Item(12, 2).Value = 162
$wsForecast.Cells.Item(12, 3).Value = -79.05050498281905
$wsForecast.Cells.Item(12, 4).Value = 392.8018728163197
$wsForecast.Cells.Item(13, 1).Value = 45361.99999999999
$wsForecast.Cells.Item(13, 2).Value = 165
$wsForecast.Cells.Item(13, 3).Value = -87.88468530577639
$wsForecast.Cells.Item(13, 4).Value = 394.6464958762446
$wsForecast.Cells.Item(14, 1).Value = 45375.99999999999
$wsForecast.Cells.Item(14, 2).Value = 169
$wsForecast.Cells.Item(14, 3).Value = -63.46155227926863
$wsForecast.Cells.Item(14, 4).Value = 397.7713359370856
$wsForecast.Cells.Item(15, 1).Value = 45382.99999999999
$wsForecast.Cells.Item(15, 2).Value = 171
$wsForecast.Cells.Item(15, 3).Value = -71.77285065858499
$wsForecast.Cells.Item(15, 4).Value = 408.5157547979476
$wsForecast.Cells.Item(16, 1).Value = 45389.99999999999
$wsForecast.Cells.Item(16, 2).Value = 173
$wsForecast.Cells.Item(16, 3).Value = -63.18214371767017
$wsForecast.Cells.Item(16, 4).Value = 391.7343116839273
$wsForecast.Cells.Item(17, 1).Value = 45403.99999999999
$wsForecast.Cells.Item(17, 2).Value = 176
$wsForecast.Cells.Item(17, 3).Value = -67.29104913802419
$wsForecast.Cells.Item(17, 4).Value = 413.6357036877841
$wsForecast.Cells.Item(18, 1).Value = 45410.99999999999
$wsForecast.Cells.Item(18, 2).Value = 178
$wsForecast.Cells.Item(18, 3).Value = -61.91276531189516
$wsForecast.Cells.Item(18, 4).Value = 414.0706513125234
$wsForecast.Cells.Item(19, 1).Value = 45417.99999999999
$wsForecast.Cells.Item(19, 2).Value = 180
$wsForecast.Cells.Item(19, 3).Value = -65.9060388082161
$wsForecast.Cells.Item(19, 4).Value = 402.4265612986914
$wsForecast.Cells.Item(20, 1).Value = 45424.99999999999
$wsForecast.Cells.Item(20, 2).Value = 182
$wsForecast.Cells.Item(20, 3).Value = -74.10903256318329
$wsForecast.Cells.Item(20, 4).Value = 419.2331720801949
$wsForecast.Cells.Item(21, 1).Value = 45431.99999999999
$wsForecast.Cells.Item(21, 2).Value = 184
$wsForecast.Cells.Item(21, 3).Value = -67.39699448627022
$wsForecast.Cells.Item(21, 4).Value = 430.0144870990749
$wsForecast.Cells.Item(22, 1).Value = 45438.99999999999
$wsForecast.Cells.Item(22, 2).Value = 186
$wsForecast.Cells.Item(22, 3).Value = -43.17916795556471
$wsForecast.Cells.Item(22, 4).Value = 433.085831325037
$wsForecast.Cells.Item(23, 1).Value = 45445.99999999999
$wsForecast.Cells.Item(23, 2).Value = 187
$wsForecast.Cells.Item(23, 3).Value = -40.58572430211465
$wsForecast.Cells.Item(23, 4).Value = 438.2461411051747
$wsForecast.Cells.Item(24, 1).Value = 45459.99999999999
$wsForecast.Cells.Item(24, 2).Value = 191
$wsForecast.Cells.Item(24, 3).Value = -49.48735633703546
$wsForecast.Cells.Item(24, 4).Value = 417.6949628322757
$wsForecast.Cells.Item(25, 1).Value = 45466.99999999999
$wsForecast.Cells.Item(25, 2).Value = 193
$wsForecast.Cells.Item(25, 3).Value = -32.39537984882926
$wsForecast.Cells.Item(25, 4).Value = 430.3868395745391
$wsForecast.Cells.Item(26, 1).Value = 45529.99999999999
$wsForecast.Cells.Item(26, 2).Value = 209
$wsForecast.Cells.Item(26, 3).Value = -18.02534135720999
$wsForecast.Cells.Item(26, 4).Value = 457.70047436465
$wsForecast.Cells.Item(27, 1).Value = 45543.99999999999
$wsForecast.Cells.Item(27, 2).Value = 213
$wsForecast.Cells.Item(27, 3).Value = -22.20511858296272
$wsForecast.Cells.Item(27, 4).Value = 452.4356258825997
$wsForecast.Cells.Item(28, 1).Value = 45564.99999999999
$wsForecast.Cells.Item(28, 2).Value = 219
$wsForecast.Cells.Item(28, 3).Value = -17.95021595527567
$wsForecast.Cells.Item(28, 4).Value = 456.1548511092986
$wsForecast.Cells.Item(29, 1).Value = 45571.99999999999
$wsForecast.Cells.Item(29, 2).Value = 220
$wsForecast.Cells.Item(29, 3).Value = -4.976154570734508
$wsForecast.Cells.Item(29, 4).Value = 455.8965514191792
$wsForecast.Cells.Item(30, 1).Value = 45578.99999999999
$wsForecast.Cells.Item(30, 2).Value = 222
$wsForecast.Cells.Item(30, 3).Value = -15.92455916911469
$wsForecast.Cells.Item(30, 4).Value = 466.1938272669167
$wsForecast.Cells.Item(31, 1).Value = 45585.99999999999
$wsForecast.Cells.Item(31, 2).Value = 224
$wsForecast.Cells.Item(31, 3).Value = -11.48455789865452
$wsForecast.Cells.Item(31, 4).Value = 460.2583401123092
$wsForecast.Cells.Item(32, 1).Value = 45599.99999999999
$wsForecast.Cells.Item(32, 2).Value = 228
$wsForecast.Cells.Item(32, 3).Value = -5.122084022273569
$wsForecast.Cells.Item(32, 4).Value = 467.7862576924372
$wsForecast.Cells.Item(33, 1).Value = 45613.99999999999
$wsForecast.Cells.Item(33, 2).Value = 231
$wsForecast.Cells.Item(33, 3).Value = 1.474116961306918
$wsForecast.Cells.Item(33, 4).Value = 463.6660751124881
$wsForecast.Cells.Item(34, 1).Value = 45620.99999999999
$wsForecast.Cells.Item(34, 2).Value = 233
$wsForecast.Cells.Item(34, 3).Value = 15.83786834965761
$wsForecast.Cells.Item(34, 4).Value = 463.7123110209289
$wsForecast.Cells.Item(35, 1).Value = 45627.99999999999
$wsForecast.Cells.Item(35, 2).Value = 235
$wsForecast.Cells.Item(35, 3).Value = -15.68109757442497
$wsForecast.Cells.Item(35, 4).Value = 467.8802183103493
$wsForecast.Cells.Item(36, 1).Value = 45634.99999999999
$wsForecast.Cells.Item(36, 2).Value = 237
$wsForecast.Cells.Item(36, 3).Value = -3.77689777825015
$wsForecast.Cells.Item(36, 4).Value = 458.9305612796688
$wsForecast.Cells.Item(37, 1).Value = 45641.99999999999
$wsForecast.Cells.Item(37, 2).Value = 239
$wsForecast.Cells.Item(37, 3).Value = -0.9959794353865125
$wsForecast.Cells.Item(37, 4).Value = 455.4752846838742
$wsForecast.Cells.Item(38, 1).Value = 45648.99999999999
$wsForecast.Cells.Item(38, 2).Value = 240
$wsForecast.Cells.Item(38, 3).Value = 2.363671124916902
$wsForecast.Cells.Item(38, 4).Value = 489.7894263950822
$wsForecast.Cells.Item(39, 1).Value = 45655.99999999999
$wsForecast.Cells.Item(39, 2).Value = 242
$wsForecast.Cells.Item(39, 3).Value = 20.09701823256697
$wsForecast.Cells.Item(39, 4).Value = 476.795856588067
$wsForecast.Cells.Item(40, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(40, 2).Value = 244
$wsForecast.Cells.Item(40, 3).Value = 5.359984753951346
$wsForecast.Cells.Item(40, 4).Value = 477.9282947275392
$wsForecast.Cells.Item(41, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(41, 2).Value = 246
$wsForecast.Cells.Item(41, 3).Value = 9.378107761443168
$wsForecast.Cells.Item(41, 4).Value = 469.2981616763078

# --- Apply styles matching other sheets (bold/border header, date format col A) ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A41").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Output "Edit applied successfully"
